$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D3").Value = "2016-02-24 09:01:03"
$wsZh.Range("G3").Value = "2016-02-24 09:02:06"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D3").Value = "2016-02-24 09:01:18"
$wsDe.Range("G3").Value = "2016-02-24 09:02:30"
